$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "THE POWER OF YOU"
$ws.Range("C2").Value = "['Somsbodvfamou ruccetful', 'Timg', 'Thi', 'blackbelt', 'Oere Mucy Cenb', 'THE POWER OF YOU']"
$ws.Range("D2").Value = "['Forbes', 'U.N']"
